# The bemonsteringsinstrument conceptscheme moved from the "omgeving" data
# domain to the "bodemenondergrond" data domain. Update every cell that
# still references the old conceptscheme URL (inScheme / topConceptOf
# columns, plus the ConceptScheme's own self-referencing id cell) so that
# they point at the new domain instead.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldUrl = "https://data.omgeving.vlaanderen.be/id/conceptscheme/bemonsteringsinstrument"
$newUrl = "https://data.bodemenondergrond.vlaanderen.be/id/conceptscheme/bemonsteringsinstrument"

$used = $ws.UsedRange
[void]$used.Replace($oldUrl, $newUrl)
